# Auto-generated edit script: apply crypto price/volume refresh + two row swaps
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.322.74"
$ws.Range("E2").Value = "  -1.19%  "
$ws.Range("D3").Value = "2.369.10"
$ws.Range("E3").Value = "  -2.20%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.86%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.636"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("E9").Value = "  -3.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.94"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0919"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.49"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.17%  "
$ws.Range("E13").Value = "  +0.02%  "
$ws.Range("E14").Value = "  -4.05%  "
$ws.Range("D15").Value = "2.728.65"
$ws.Range("E15").Value = "  -1.59%  "
$ws.Range("E16").Value = "  -3.87%  "
$ws.Range("D17").Value = "2.368.93"
$ws.Range("E17").Value = "  -1.94%  "
$ws.Range("D18").Value = "45.256.19"
$ws.Range("E18").Value = "  -1.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "15.58"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +15.11%  "
$ws.Range("E20").Value = "  -4.73%  "
$ws.Range("E21").Value = "  -2.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.59"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.52%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "263.84"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.39%  "
$ws.Range("E25").Value = "  -1.20%  "
$ws.Range("E26").Value = "  +0.16%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.29%  "
$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.51"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.62%  "
$ws.Range("E29").Value = "  -1.74%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.41"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0949"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "37.34"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.99%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "168.75"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.88"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.87%  "
$ws.Range("E35").Value = "  -0.77%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.117"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.72"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.35%  "
$ws.Range("E38").Value = "  -1.91%  "
$ws.Range("E39").Value = "  +7.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.02"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.82%  "
$ws.Range("E41").Value = "  -3.66%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.23"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "70.13"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.71%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "1.872.27"
$ws.Range("E44").Value = "  +13.20%  "
$ws.Range("B45").Value = "Celestia"
$ws.Range("C45").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.98"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.228"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.50%  "
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.98"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.64%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "84.81"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.63%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "112.20"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.25%  "
$ws.Range("E51").Value = "  -2.34%  "
